$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "ДОР.ОП"
$ws.Range("B14").Value = "Отдел НИР"

$ws.Range("B15").Select()
